$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2..49: update id (B), speaker_variant (C), and clear is_prefered (D)
$ws.Cells.Item(2, 2).Value = "#verstant"
$ws.Cells.Item(2, 3).Value = "Verstant"
$ws.Cells.Item(2, 4).ClearContents()
$ws.Cells.Item(3, 2).Value = "#tyter"
$ws.Cells.Item(3, 3).Value = "Tyter"
$ws.Cells.Item(3, 4).ClearContents()
$ws.Cells.Item(4, 2).Value = "#kluysen"
$ws.Cells.Item(4, 3).Value = "Kluysen"
$ws.Cells.Item(4, 4).ClearContents()
$ws.Cells.Item(5, 2).Value = "#octavio"
$ws.Cells.Item(5, 3).Value = "Octavio"
$ws.Cells.Item(5, 4).ClearContents()
$ws.Cells.Item(6, 2).Value = "#vville"
$ws.Cells.Item(6, 3).Value = "VVille"
$ws.Cells.Item(6, 4).ClearContents()
$ws.Cells.Item(7, 2).Value = "#ottavio"
$ws.Cells.Item(7, 3).Value = "Ottavio"
$ws.Cells.Item(7, 4).ClearContents()
$ws.Cells.Item(8, 2).Value = "#theoph"
$ws.Cells.Item(8, 3).Value = "Theoph"
$ws.Cells.Item(8, 4).ClearContents()
$ws.Cells.Item(9, 2).Value = "#af-sijnde"
$ws.Cells.Item(9, 3).Value = "Af-sijnde"
$ws.Cells.Item(9, 4).ClearContents()
$ws.Cells.Item(10, 2).Value = "#nietg"
$ws.Cells.Item(10, 3).Value = "Nietg"
$ws.Cells.Item(10, 4).ClearContents()
$ws.Cells.Item(11, 2).Value = "#guyd"
$ws.Cells.Item(11, 3).Value = "Guyd"
$ws.Cells.Item(11, 4).ClearContents()
$ws.Cells.Item(12, 2).Value = "#gualdr"
$ws.Cells.Item(12, 3).Value = "Gualdr"
$ws.Cells.Item(12, 4).ClearContents()
$ws.Cells.Item(13, 2).Value = "#nietgel"
$ws.Cells.Item(13, 3).Value = "Nietgel"
$ws.Cells.Item(13, 4).ClearContents()
$ws.Cells.Item(14, 2).Value = "#cupido"
$ws.Cells.Item(14, 3).Value = "Cupido"
$ws.Cells.Item(14, 4).ClearContents()
$ws.Cells.Item(15, 2).Value = "#arae"
$ws.Cells.Item(15, 3).Value = "Arae"
$ws.Cells.Item(15, 4).ClearContents()
$ws.Cells.Item(16, 2).Value = "#'tspels-inhout"
$ws.Cells.Item(16, 3).Value = "''Tspels inhout"
$ws.Cells.Item(16, 4).ClearContents()
$ws.Cells.Item(17, 2).Value = "#brand"
$ws.Cells.Item(17, 3).Value = "Brand"
$ws.Cells.Item(17, 4).ClearContents()
$ws.Cells.Item(18, 2).Value = "#vvil"
$ws.Cells.Item(18, 3).Value = "VVil"
$ws.Cells.Item(18, 4).ClearContents()
$ws.Cells.Item(19, 2).Value = "#guydeon"
$ws.Cells.Item(19, 3).Value = "Guydeon"
$ws.Cells.Item(19, 4).ClearContents()
$ws.Cells.Item(20, 2).Value = "#niet-ghe"
$ws.Cells.Item(20, 3).Value = "Niet-ghe"
$ws.Cells.Item(20, 4).ClearContents()
$ws.Cells.Item(21, 2).Value = "#laura,"
$ws.Cells.Item(21, 3).Value = "Laura,"
$ws.Cells.Item(21, 4).ClearContents()
$ws.Cells.Item(22, 2).Value = "#billinc"
$ws.Cells.Item(22, 3).Value = "Billinc"
$ws.Cells.Item(22, 4).ClearContents()
$ws.Cells.Item(23, 2).Value = "#keyser"
$ws.Cells.Item(23, 3).Value = "Keyser"
$ws.Cells.Item(23, 4).ClearContents()
$ws.Cells.Item(24, 2).Value = "#galdra"
$ws.Cells.Item(24, 3).Value = "Galdra"
$ws.Cells.Item(24, 4).ClearContents()
$ws.Cells.Item(25, 2).Value = "#kluysenaer,"
$ws.Cells.Item(25, 3).Value = "Kluysenaer,"
$ws.Cells.Item(25, 4).ClearContents()
$ws.Cells.Item(26, 2).Value = "#af-zijnde"
$ws.Cells.Item(26, 3).Value = "Af-zijnde"
$ws.Cells.Item(26, 4).ClearContents()
$ws.Cells.Item(27, 2).Value = "#branden"
$ws.Cells.Item(27, 3).Value = "Branden"
$ws.Cells.Item(27, 4).ClearContents()
$ws.Cells.Item(28, 2).Value = "#gald"
$ws.Cells.Item(28, 3).Value = "Gald"
$ws.Cells.Item(28, 4).ClearContents()
$ws.Cells.Item(29, 2).Value = "#theophe"
$ws.Cells.Item(29, 3).Value = "Theophe"
$ws.Cells.Item(29, 4).ClearContents()
$ws.Cells.Item(30, 2).Value = "#kluysena"
$ws.Cells.Item(30, 3).Value = "Kluysena"
$ws.Cells.Item(30, 4).ClearContents()
$ws.Cells.Item(31, 2).Value = "#flav"
$ws.Cells.Item(31, 3).Value = "Flav"
$ws.Cells.Item(31, 4).ClearContents()
$ws.Cells.Item(32, 2).Value = "#kluyse"
$ws.Cells.Item(32, 3).Value = "Kluyse"
$ws.Cells.Item(32, 4).ClearContents()
$ws.Cells.Item(33, 2).Value = "#galdrad"
$ws.Cells.Item(33, 3).Value = "Galdrad"
$ws.Cells.Item(33, 4).ClearContents()
$ws.Cells.Item(34, 2).Value = "#lust"
$ws.Cells.Item(34, 3).Value = "Lust"
$ws.Cells.Item(34, 4).ClearContents()
$ws.Cells.Item(35, 2).Value = "#ialourse"
$ws.Cells.Item(35, 3).Value = "Ialourse"
$ws.Cells.Item(35, 4).ClearContents()
$ws.Cells.Item(36, 2).Value = "#kluysenaer"
$ws.Cells.Item(36, 3).Value = "Kluysenaer"
$ws.Cells.Item(36, 4).ClearContents()
$ws.Cells.Item(37, 2).Value = "#galdr"
$ws.Cells.Item(37, 3).Value = "Galdr"
$ws.Cells.Item(37, 4).ClearContents()
$ws.Cells.Item(38, 2).Value = "#bedil"
$ws.Cells.Item(38, 3).Value = "Bedil"
$ws.Cells.Item(38, 4).ClearContents()
$ws.Cells.Item(39, 2).Value = "#flavio"
$ws.Cells.Item(39, 3).Value = "Flavio"
$ws.Cells.Item(39, 4).ClearContents()
$ws.Cells.Item(40, 2).Value = "#bellinc"
$ws.Cells.Item(40, 3).Value = "Bellinc"
$ws.Cells.Item(40, 4).ClearContents()
$ws.Cells.Item(41, 2).Value = "#camillo"
$ws.Cells.Item(41, 3).Value = "Camillo"
$ws.Cells.Item(41, 4).ClearContents()
$ws.Cells.Item(42, 2).Value = "#flavio,"
$ws.Cells.Item(42, 3).Value = "Flavio,"
$ws.Cells.Item(42, 4).ClearContents()
$ws.Cells.Item(43, 2).Value = "#laura"
$ws.Cells.Item(43, 3).Value = "Laura"
$ws.Cells.Item(43, 4).ClearContents()
$ws.Cells.Item(44, 2).Value = "#niet-ghel"
$ws.Cells.Item(44, 3).Value = "Niet-ghel"
$ws.Cells.Item(44, 4).ClearContents()
$ws.Cells.Item(45, 2).Value = "#araeta"
$ws.Cells.Item(45, 3).Value = "Araeta"
$ws.Cells.Item(45, 4).ClearContents()
$ws.Cells.Item(46, 2).Value = "#brandenb"
$ws.Cells.Item(46, 3).Value = "Brandenb"
$ws.Cells.Item(46, 4).ClearContents()
$ws.Cells.Item(47, 2).Value = "#bedil-al"
$ws.Cells.Item(47, 3).Value = "Bedil-al"
$ws.Cells.Item(47, 4).ClearContents()
$ws.Cells.Item(48, 2).Value = "#galdrada"
$ws.Cells.Item(48, 3).Value = "Galdrada"
$ws.Cells.Item(48, 4).ClearContents()
$ws.Cells.Item(49, 2).Value = "#niet-geliefde"
$ws.Cells.Item(49, 3).Value = "Niet-geliefde"
$ws.Cells.Item(49, 4).ClearContents()
